$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: user "fahrezi" -> "alcoy"
$ws.Range("B2").Value = "alcoy"

# Row 2: dates "2024-02-15" -> "2024-02-16"
# Row 3 (D3/E3) already holds the text "2024-02-16" as a shared string,
# so copy it across to avoid Excel re-parsing the literal as a real date.
$ws.Range("D3").Copy($ws.Range("D2"))
$ws.Range("E3").Copy($ws.Range("E2"))

# Row 3: user "alcoy" -> "fahrezi"
$ws.Range("B3").Value = "fahrezi"
